$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows got shifted down by one data point starting at row 69:
# row 69 receives a brand new observation, and every row 70..156 now holds
# what used to be one row above it (row 156 is a brand-new row holding the
# old row 155 values). Columns A,B,C,E,F,G,H,I,O,R are identical on every
# row in this sheet, so only D,J,K,L,M,N,P,Q actually change per row; row 156
# additionally needs the constant columns since it is a new row.

# Populate the brand-new row 156 with the constant columns shared by all rows
$ws.Cells.Item(156, 1).Value = 8                                            # A156
$ws.Cells.Item(156, 2).Value = "Terminal La Palmera de La Serena"           # B156
$ws.Cells.Item(156, 3).Value = "Coquimbo"                                   # C156
$ws.Cells.Item(156, 5).Value = 4                                            # E156
$ws.Cells.Item(156, 6).Value = 100112001                                    # F156
$ws.Cells.Item(156, 7).Value = "Berenjena"                                  # G156
$ws.Cells.Item(156, 8).Value = "Sin especificar"                            # H156
$ws.Cells.Item(156, 9).Value = "Primera"                                    # I156
$ws.Cells.Item(156, 15).Value = "Región de Arica y Parinacota"              # O156
$ws.Cells.Item(156, 18).Value = "Hortaliza"                                 # R156

# Make row 156 use the same date-number-format style as the other D column cells
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(155, 4).NumberFormat

# Now update the shifted-down data columns (D, J, K, L, M, N, P, Q) for rows 69-156

$ws.Cells.Item(69, 4).Value = 44789    # D69
$ws.Cells.Item(69, 10).Value = 500   # J69
$ws.Cells.Item(69, 11).Value = 10500   # K69
$ws.Cells.Item(69, 12).Value = 11000   # L69
$ws.Cells.Item(69, 13).Value = 10750   # M69
$ws.Cells.Item(69, 14).Value = "`$/caja 50 unidades"   # N69
$ws.Cells.Item(69, 16).Value = 215   # P69
$ws.Cells.Item(69, 17).Value = 50   # Q69

$ws.Cells.Item(70, 4).Value = 44568    # D70
$ws.Cells.Item(70, 10).Value = 700   # J70
$ws.Cells.Item(70, 11).Value = 8000   # K70
$ws.Cells.Item(70, 12).Value = 9000   # L70
$ws.Cells.Item(70, 13).Value = 8500   # M70
$ws.Cells.Item(70, 14).Value = "`$/caja 60 unidades"   # N70
$ws.Cells.Item(70, 16).Value = 142   # P70
$ws.Cells.Item(70, 17).Value = 60   # Q70

$ws.Cells.Item(71, 4).Value = 44746    # D71
$ws.Cells.Item(71, 10).Value = 700   # J71
$ws.Cells.Item(71, 11).Value = 9000   # K71
$ws.Cells.Item(71, 12).Value = 10000   # L71
$ws.Cells.Item(71, 13).Value = 9500   # M71
$ws.Cells.Item(71, 14).Value = "`$/caja 50 unidades"   # N71
$ws.Cells.Item(71, 16).Value = 190   # P71
$ws.Cells.Item(71, 17).Value = 50   # Q71

$ws.Cells.Item(72, 4).Value = 44258    # D72
$ws.Cells.Item(72, 10).Value = 700   # J72
$ws.Cells.Item(72, 11).Value = 8000   # K72
$ws.Cells.Item(72, 12).Value = 8500   # L72
$ws.Cells.Item(72, 13).Value = 8250   # M72
$ws.Cells.Item(72, 14).Value = "`$/caja 60 unidades"   # N72
$ws.Cells.Item(72, 16).Value = 138   # P72
$ws.Cells.Item(72, 17).Value = 60   # Q72

$ws.Cells.Item(73, 4).Value = 44323    # D73
$ws.Cells.Item(73, 10).Value = 500   # J73
$ws.Cells.Item(73, 11).Value = 8000   # K73
$ws.Cells.Item(73, 12).Value = 9000   # L73
$ws.Cells.Item(73, 13).Value = 8500   # M73
$ws.Cells.Item(73, 14).Value = "`$/caja 60 unidades"   # N73
$ws.Cells.Item(73, 16).Value = 142   # P73
$ws.Cells.Item(73, 17).Value = 60   # Q73

$ws.Cells.Item(74, 4).Value = 44477    # D74
$ws.Cells.Item(74, 10).Value = 600   # J74
$ws.Cells.Item(74, 11).Value = 8000   # K74
$ws.Cells.Item(74, 12).Value = 9000   # L74
$ws.Cells.Item(74, 13).Value = 8500   # M74
$ws.Cells.Item(74, 14).Value = "`$/caja 60 unidades"   # N74
$ws.Cells.Item(74, 16).Value = 142   # P74
$ws.Cells.Item(74, 17).Value = 60   # Q74

$ws.Cells.Item(75, 4).Value = 44643    # D75
$ws.Cells.Item(75, 10).Value = 560   # J75
$ws.Cells.Item(75, 11).Value = 8000   # K75
$ws.Cells.Item(75, 12).Value = 9000   # L75
$ws.Cells.Item(75, 13).Value = 8500   # M75
$ws.Cells.Item(75, 14).Value = "`$/caja 50 unidades"   # N75
$ws.Cells.Item(75, 16).Value = 170   # P75
$ws.Cells.Item(75, 17).Value = 50   # Q75

$ws.Cells.Item(76, 4).Value = 44403    # D76
$ws.Cells.Item(76, 10).Value = 760   # J76
$ws.Cells.Item(76, 11).Value = 12000   # K76
$ws.Cells.Item(76, 12).Value = 12500   # L76
$ws.Cells.Item(76, 13).Value = 12250   # M76
$ws.Cells.Item(76, 14).Value = "`$/caja 60 unidades"   # N76
$ws.Cells.Item(76, 16).Value = 204   # P76
$ws.Cells.Item(76, 17).Value = 60   # Q76

$ws.Cells.Item(77, 4).Value = 44407    # D77
$ws.Cells.Item(77, 10).Value = 600   # J77
$ws.Cells.Item(77, 11).Value = 12500   # K77
$ws.Cells.Item(77, 12).Value = 13000   # L77
$ws.Cells.Item(77, 13).Value = 12750   # M77
$ws.Cells.Item(77, 14).Value = "`$/caja 60 unidades"   # N77
$ws.Cells.Item(77, 16).Value = 212   # P77
$ws.Cells.Item(77, 17).Value = 60   # Q77

$ws.Cells.Item(78, 4).Value = 44657    # D78
$ws.Cells.Item(78, 10).Value = 2000   # J78
$ws.Cells.Item(78, 11).Value = 8500   # K78
$ws.Cells.Item(78, 12).Value = 9000   # L78
$ws.Cells.Item(78, 13).Value = 8750   # M78
$ws.Cells.Item(78, 14).Value = "`$/caja 50 unidades"   # N78
$ws.Cells.Item(78, 16).Value = 175   # P78
$ws.Cells.Item(78, 17).Value = 50   # Q78

$ws.Cells.Item(79, 4).Value = 44721    # D79
$ws.Cells.Item(79, 10).Value = 520   # J79
$ws.Cells.Item(79, 11).Value = 8000   # K79
$ws.Cells.Item(79, 12).Value = 9000   # L79
$ws.Cells.Item(79, 13).Value = 8500   # M79
$ws.Cells.Item(79, 14).Value = "`$/caja 50 unidades"   # N79
$ws.Cells.Item(79, 16).Value = 170   # P79
$ws.Cells.Item(79, 17).Value = 50   # Q79

$ws.Cells.Item(80, 4).Value = 44554    # D80
$ws.Cells.Item(80, 10).Value = 500   # J80
$ws.Cells.Item(80, 11).Value = 10000   # K80
$ws.Cells.Item(80, 12).Value = 11000   # L80
$ws.Cells.Item(80, 13).Value = 10500   # M80
$ws.Cells.Item(80, 14).Value = "`$/caja 60 unidades"   # N80
$ws.Cells.Item(80, 16).Value = 175   # P80
$ws.Cells.Item(80, 17).Value = 60   # Q80

$ws.Cells.Item(81, 4).Value = 44377    # D81
$ws.Cells.Item(81, 10).Value = 600   # J81
$ws.Cells.Item(81, 11).Value = 12000   # K81
$ws.Cells.Item(81, 12).Value = 13000   # L81
$ws.Cells.Item(81, 13).Value = 12500   # M81
$ws.Cells.Item(81, 14).Value = "`$/caja 60 unidades"   # N81
$ws.Cells.Item(81, 16).Value = 208   # P81
$ws.Cells.Item(81, 17).Value = 60   # Q81

$ws.Cells.Item(82, 4).Value = 44587    # D82
$ws.Cells.Item(82, 10).Value = 520   # J82
$ws.Cells.Item(82, 11).Value = 8000   # K82
$ws.Cells.Item(82, 12).Value = 9000   # L82
$ws.Cells.Item(82, 13).Value = 8500   # M82
$ws.Cells.Item(82, 14).Value = "`$/caja 50 unidades"   # N82
$ws.Cells.Item(82, 16).Value = 170   # P82
$ws.Cells.Item(82, 17).Value = 50   # Q82

$ws.Cells.Item(83, 4).Value = 44767    # D83
$ws.Cells.Item(83, 10).Value = 600   # J83
$ws.Cells.Item(83, 11).Value = 10000   # K83
$ws.Cells.Item(83, 12).Value = 11000   # L83
$ws.Cells.Item(83, 13).Value = 10500   # M83
$ws.Cells.Item(83, 14).Value = "`$/caja 50 unidades"   # N83
$ws.Cells.Item(83, 16).Value = 210   # P83
$ws.Cells.Item(83, 17).Value = 50   # Q83

$ws.Cells.Item(84, 4).Value = 44503    # D84
$ws.Cells.Item(84, 10).Value = 600   # J84
$ws.Cells.Item(84, 11).Value = 8000   # K84
$ws.Cells.Item(84, 12).Value = 8500   # L84
$ws.Cells.Item(84, 13).Value = 8250   # M84
$ws.Cells.Item(84, 14).Value = "`$/caja 60 unidades"   # N84
$ws.Cells.Item(84, 16).Value = 138   # P84
$ws.Cells.Item(84, 17).Value = 60   # Q84

$ws.Cells.Item(85, 4).Value = 44307    # D85
$ws.Cells.Item(85, 10).Value = 600   # J85
$ws.Cells.Item(85, 11).Value = 8000   # K85
$ws.Cells.Item(85, 12).Value = 9000   # L85
$ws.Cells.Item(85, 13).Value = 8500   # M85
$ws.Cells.Item(85, 14).Value = "`$/caja 60 unidades"   # N85
$ws.Cells.Item(85, 16).Value = 142   # P85
$ws.Cells.Item(85, 17).Value = 60   # Q85

$ws.Cells.Item(86, 4).Value = 44498    # D86
$ws.Cells.Item(86, 10).Value = 560   # J86
$ws.Cells.Item(86, 11).Value = 8000   # K86
$ws.Cells.Item(86, 12).Value = 9000   # L86
$ws.Cells.Item(86, 13).Value = 8500   # M86
$ws.Cells.Item(86, 14).Value = "`$/caja 60 unidades"   # N86
$ws.Cells.Item(86, 16).Value = 142   # P86
$ws.Cells.Item(86, 17).Value = 60   # Q86

$ws.Cells.Item(87, 4).Value = 44466    # D87
$ws.Cells.Item(87, 10).Value = 600   # J87
$ws.Cells.Item(87, 11).Value = 8500   # K87
$ws.Cells.Item(87, 12).Value = 9000   # L87
$ws.Cells.Item(87, 13).Value = 8750   # M87
$ws.Cells.Item(87, 14).Value = "`$/caja 60 unidades"   # N87
$ws.Cells.Item(87, 16).Value = 146   # P87
$ws.Cells.Item(87, 17).Value = 60   # Q87

$ws.Cells.Item(88, 4).Value = 44427    # D88
$ws.Cells.Item(88, 10).Value = 560   # J88
$ws.Cells.Item(88, 11).Value = 12000   # K88
$ws.Cells.Item(88, 12).Value = 13000   # L88
$ws.Cells.Item(88, 13).Value = 12500   # M88
$ws.Cells.Item(88, 14).Value = "`$/caja 60 unidades"   # N88
$ws.Cells.Item(88, 16).Value = 208   # P88
$ws.Cells.Item(88, 17).Value = 60   # Q88

$ws.Cells.Item(89, 4).Value = 44349    # D89
$ws.Cells.Item(89, 10).Value = 600   # J89
$ws.Cells.Item(89, 11).Value = 12000   # K89
$ws.Cells.Item(89, 12).Value = 12500   # L89
$ws.Cells.Item(89, 13).Value = 12250   # M89
$ws.Cells.Item(89, 14).Value = "`$/caja 60 unidades"   # N89
$ws.Cells.Item(89, 16).Value = 204   # P89
$ws.Cells.Item(89, 17).Value = 60   # Q89

$ws.Cells.Item(90, 4).Value = 44545    # D90
$ws.Cells.Item(90, 10).Value = 540   # J90
$ws.Cells.Item(90, 11).Value = 10000   # K90
$ws.Cells.Item(90, 12).Value = 11000   # L90
$ws.Cells.Item(90, 13).Value = 10500   # M90
$ws.Cells.Item(90, 14).Value = "`$/caja 60 unidades"   # N90
$ws.Cells.Item(90, 16).Value = 175   # P90
$ws.Cells.Item(90, 17).Value = 60   # Q90

$ws.Cells.Item(91, 4).Value = 44326    # D91
$ws.Cells.Item(91, 10).Value = 500   # J91
$ws.Cells.Item(91, 11).Value = 9000   # K91
$ws.Cells.Item(91, 12).Value = 10000   # L91
$ws.Cells.Item(91, 13).Value = 9500   # M91
$ws.Cells.Item(91, 14).Value = "`$/caja 60 unidades"   # N91
$ws.Cells.Item(91, 16).Value = 158   # P91
$ws.Cells.Item(91, 17).Value = 60   # Q91

$ws.Cells.Item(92, 4).Value = 44384    # D92
$ws.Cells.Item(92, 10).Value = 600   # J92
$ws.Cells.Item(92, 11).Value = 11000   # K92
$ws.Cells.Item(92, 12).Value = 12000   # L92
$ws.Cells.Item(92, 13).Value = 11500   # M92
$ws.Cells.Item(92, 14).Value = "`$/caja 60 unidades"   # N92
$ws.Cells.Item(92, 16).Value = 192   # P92
$ws.Cells.Item(92, 17).Value = 60   # Q92

$ws.Cells.Item(93, 4).Value = 44314    # D93
$ws.Cells.Item(93, 10).Value = 560   # J93
$ws.Cells.Item(93, 11).Value = 8000   # K93
$ws.Cells.Item(93, 12).Value = 9000   # L93
$ws.Cells.Item(93, 13).Value = 8500   # M93
$ws.Cells.Item(93, 14).Value = "`$/caja 60 unidades"   # N93
$ws.Cells.Item(93, 16).Value = 142   # P93
$ws.Cells.Item(93, 17).Value = 60   # Q93

$ws.Cells.Item(94, 4).Value = 44496    # D94
$ws.Cells.Item(94, 10).Value = 520   # J94
$ws.Cells.Item(94, 11).Value = 8000   # K94
$ws.Cells.Item(94, 12).Value = 9000   # L94
$ws.Cells.Item(94, 13).Value = 8500   # M94
$ws.Cells.Item(94, 14).Value = "`$/caja 60 unidades"   # N94
$ws.Cells.Item(94, 16).Value = 142   # P94
$ws.Cells.Item(94, 17).Value = 60   # Q94

$ws.Cells.Item(95, 4).Value = 44293    # D95
$ws.Cells.Item(95, 10).Value = 600   # J95
$ws.Cells.Item(95, 11).Value = 8000   # K95
$ws.Cells.Item(95, 12).Value = 9000   # L95
$ws.Cells.Item(95, 13).Value = 8500   # M95
$ws.Cells.Item(95, 14).Value = "`$/caja 60 unidades"   # N95
$ws.Cells.Item(95, 16).Value = 142   # P95
$ws.Cells.Item(95, 17).Value = 60   # Q95

$ws.Cells.Item(96, 4).Value = 44357    # D96
$ws.Cells.Item(96, 10).Value = 520   # J96
$ws.Cells.Item(96, 11).Value = 12000   # K96
$ws.Cells.Item(96, 12).Value = 12500   # L96
$ws.Cells.Item(96, 13).Value = 12250   # M96
$ws.Cells.Item(96, 14).Value = "`$/caja 60 unidades"   # N96
$ws.Cells.Item(96, 16).Value = 204   # P96
$ws.Cells.Item(96, 17).Value = 60   # Q96

$ws.Cells.Item(97, 4).Value = 44524    # D97
$ws.Cells.Item(97, 10).Value = 540   # J97
$ws.Cells.Item(97, 11).Value = 9000   # K97
$ws.Cells.Item(97, 12).Value = 10000   # L97
$ws.Cells.Item(97, 13).Value = 9500   # M97
$ws.Cells.Item(97, 14).Value = "`$/caja 60 unidades"   # N97
$ws.Cells.Item(97, 16).Value = 158   # P97
$ws.Cells.Item(97, 17).Value = 60   # Q97

$ws.Cells.Item(98, 4).Value = 44321    # D98
$ws.Cells.Item(98, 10).Value = 600   # J98
$ws.Cells.Item(98, 11).Value = 8000   # K98
$ws.Cells.Item(98, 12).Value = 9000   # L98
$ws.Cells.Item(98, 13).Value = 8500   # M98
$ws.Cells.Item(98, 14).Value = "`$/caja 60 unidades"   # N98
$ws.Cells.Item(98, 16).Value = 142   # P98
$ws.Cells.Item(98, 17).Value = 60   # Q98

$ws.Cells.Item(99, 4).Value = 44414    # D99
$ws.Cells.Item(99, 10).Value = 600   # J99
$ws.Cells.Item(99, 11).Value = 12500   # K99
$ws.Cells.Item(99, 12).Value = 13000   # L99
$ws.Cells.Item(99, 13).Value = 12750   # M99
$ws.Cells.Item(99, 14).Value = "`$/caja 60 unidades"   # N99
$ws.Cells.Item(99, 16).Value = 212   # P99
$ws.Cells.Item(99, 17).Value = 60   # Q99

$ws.Cells.Item(100, 4).Value = 44664    # D100
$ws.Cells.Item(100, 10).Value = 520   # J100
$ws.Cells.Item(100, 11).Value = 8000   # K100
$ws.Cells.Item(100, 12).Value = 9000   # L100
$ws.Cells.Item(100, 13).Value = 8500   # M100
$ws.Cells.Item(100, 14).Value = "`$/caja 50 unidades"   # N100
$ws.Cells.Item(100, 16).Value = 170   # P100
$ws.Cells.Item(100, 17).Value = 50   # Q100

$ws.Cells.Item(101, 4).Value = 44512    # D101
$ws.Cells.Item(101, 10).Value = 600   # J101
$ws.Cells.Item(101, 11).Value = 8000   # K101
$ws.Cells.Item(101, 12).Value = 9000   # L101
$ws.Cells.Item(101, 13).Value = 8500   # M101
$ws.Cells.Item(101, 14).Value = "`$/caja 60 unidades"   # N101
$ws.Cells.Item(101, 16).Value = 142   # P101
$ws.Cells.Item(101, 17).Value = 60   # Q101

$ws.Cells.Item(102, 4).Value = 44715    # D102
$ws.Cells.Item(102, 10).Value = 1000   # J102
$ws.Cells.Item(102, 11).Value = 8000   # K102
$ws.Cells.Item(102, 12).Value = 9000   # L102
$ws.Cells.Item(102, 13).Value = 8500   # M102
$ws.Cells.Item(102, 14).Value = "`$/caja 50 unidades"   # N102
$ws.Cells.Item(102, 16).Value = 170   # P102
$ws.Cells.Item(102, 17).Value = 50   # Q102

$ws.Cells.Item(103, 4).Value = 44582    # D103
$ws.Cells.Item(103, 10).Value = 600   # J103
$ws.Cells.Item(103, 11).Value = 8500   # K103
$ws.Cells.Item(103, 12).Value = 9000   # L103
$ws.Cells.Item(103, 13).Value = 8750   # M103
$ws.Cells.Item(103, 14).Value = "`$/caja 50 unidades"   # N103
$ws.Cells.Item(103, 16).Value = 175   # P103
$ws.Cells.Item(103, 17).Value = 50   # Q103

$ws.Cells.Item(104, 4).Value = 44300    # D104
$ws.Cells.Item(104, 10).Value = 600   # J104
$ws.Cells.Item(104, 11).Value = 8000   # K104
$ws.Cells.Item(104, 12).Value = 8500   # L104
$ws.Cells.Item(104, 13).Value = 8250   # M104
$ws.Cells.Item(104, 14).Value = "`$/caja 60 unidades"   # N104
$ws.Cells.Item(104, 16).Value = 138   # P104
$ws.Cells.Item(104, 17).Value = 60   # Q104

$ws.Cells.Item(105, 4).Value = 44656    # D105
$ws.Cells.Item(105, 10).Value = 400   # J105
$ws.Cells.Item(105, 11).Value = 8000   # K105
$ws.Cells.Item(105, 12).Value = 9000   # L105
$ws.Cells.Item(105, 13).Value = 8500   # M105
$ws.Cells.Item(105, 14).Value = "`$/caja 50 unidades"   # N105
$ws.Cells.Item(105, 16).Value = 170   # P105
$ws.Cells.Item(105, 17).Value = 50   # Q105

$ws.Cells.Item(106, 4).Value = 44519    # D106
$ws.Cells.Item(106, 10).Value = 560   # J106
$ws.Cells.Item(106, 11).Value = 8000   # K106
$ws.Cells.Item(106, 12).Value = 8500   # L106
$ws.Cells.Item(106, 13).Value = 8250   # M106
$ws.Cells.Item(106, 14).Value = "`$/caja 60 unidades"   # N106
$ws.Cells.Item(106, 16).Value = 138   # P106
$ws.Cells.Item(106, 17).Value = 60   # Q106

$ws.Cells.Item(107, 4).Value = 44676    # D107
$ws.Cells.Item(107, 10).Value = 2600   # J107
$ws.Cells.Item(107, 11).Value = 9000   # K107
$ws.Cells.Item(107, 12).Value = 10000   # L107
$ws.Cells.Item(107, 13).Value = 9500   # M107
$ws.Cells.Item(107, 14).Value = "`$/caja 50 unidades"   # N107
$ws.Cells.Item(107, 16).Value = 190   # P107
$ws.Cells.Item(107, 17).Value = 50   # Q107

$ws.Cells.Item(108, 4).Value = 44508    # D108
$ws.Cells.Item(108, 10).Value = 520   # J108
$ws.Cells.Item(108, 11).Value = 8000   # K108
$ws.Cells.Item(108, 12).Value = 9000   # L108
$ws.Cells.Item(108, 13).Value = 8500   # M108
$ws.Cells.Item(108, 14).Value = "`$/caja 60 unidades"   # N108
$ws.Cells.Item(108, 16).Value = 142   # P108
$ws.Cells.Item(108, 17).Value = 60   # Q108

$ws.Cells.Item(109, 4).Value = 44335    # D109
$ws.Cells.Item(109, 10).Value = 600   # J109
$ws.Cells.Item(109, 11).Value = 12000   # K109
$ws.Cells.Item(109, 12).Value = 13000   # L109
$ws.Cells.Item(109, 13).Value = 12500   # M109
$ws.Cells.Item(109, 14).Value = "`$/caja 60 unidades"   # N109
$ws.Cells.Item(109, 16).Value = 208   # P109
$ws.Cells.Item(109, 17).Value = 60   # Q109

$ws.Cells.Item(110, 4).Value = 44315    # D110
$ws.Cells.Item(110, 10).Value = 440   # J110
$ws.Cells.Item(110, 11).Value = 8000   # K110
$ws.Cells.Item(110, 12).Value = 9000   # L110
$ws.Cells.Item(110, 13).Value = 8500   # M110
$ws.Cells.Item(110, 14).Value = "`$/caja 60 unidades"   # N110
$ws.Cells.Item(110, 16).Value = 142   # P110
$ws.Cells.Item(110, 17).Value = 60   # Q110

$ws.Cells.Item(111, 4).Value = 44386    # D111
$ws.Cells.Item(111, 10).Value = 560   # J111
$ws.Cells.Item(111, 11).Value = 11000   # K111
$ws.Cells.Item(111, 12).Value = 12000   # L111
$ws.Cells.Item(111, 13).Value = 11500   # M111
$ws.Cells.Item(111, 14).Value = "`$/caja 60 unidades"   # N111
$ws.Cells.Item(111, 16).Value = 192   # P111
$ws.Cells.Item(111, 17).Value = 60   # Q111

$ws.Cells.Item(112, 4).Value = 44327    # D112
$ws.Cells.Item(112, 10).Value = 500   # J112
$ws.Cells.Item(112, 11).Value = 9000   # K112
$ws.Cells.Item(112, 12).Value = 10000   # L112
$ws.Cells.Item(112, 13).Value = 9500   # M112
$ws.Cells.Item(112, 14).Value = "`$/caja 60 unidades"   # N112
$ws.Cells.Item(112, 16).Value = 158   # P112
$ws.Cells.Item(112, 17).Value = 60   # Q112

$ws.Cells.Item(113, 4).Value = 44316    # D113
$ws.Cells.Item(113, 10).Value = 520   # J113
$ws.Cells.Item(113, 11).Value = 8000   # K113
$ws.Cells.Item(113, 12).Value = 9000   # L113
$ws.Cells.Item(113, 13).Value = 8500   # M113
$ws.Cells.Item(113, 14).Value = "`$/caja 60 unidades"   # N113
$ws.Cells.Item(113, 16).Value = 142   # P113
$ws.Cells.Item(113, 17).Value = 60   # Q113

$ws.Cells.Item(114, 4).Value = 44586    # D114
$ws.Cells.Item(114, 10).Value = 600   # J114
$ws.Cells.Item(114, 11).Value = 8000   # K114
$ws.Cells.Item(114, 12).Value = 9000   # L114
$ws.Cells.Item(114, 13).Value = 8500   # M114
$ws.Cells.Item(114, 14).Value = "`$/caja 50 unidades"   # N114
$ws.Cells.Item(114, 16).Value = 170   # P114
$ws.Cells.Item(114, 17).Value = 50   # Q114

$ws.Cells.Item(115, 4).Value = 44320    # D115
$ws.Cells.Item(115, 10).Value = 520   # J115
$ws.Cells.Item(115, 11).Value = 8000   # K115
$ws.Cells.Item(115, 12).Value = 9000   # L115
$ws.Cells.Item(115, 13).Value = 8500   # M115
$ws.Cells.Item(115, 14).Value = "`$/caja 60 unidades"   # N115
$ws.Cells.Item(115, 16).Value = 142   # P115
$ws.Cells.Item(115, 17).Value = 60   # Q115

$ws.Cells.Item(116, 4).Value = 44566    # D116
$ws.Cells.Item(116, 10).Value = 540   # J116
$ws.Cells.Item(116, 11).Value = 9500   # K116
$ws.Cells.Item(116, 12).Value = 10000   # L116
$ws.Cells.Item(116, 13).Value = 9750   # M116
$ws.Cells.Item(116, 14).Value = "`$/caja 60 unidades"   # N116
$ws.Cells.Item(116, 16).Value = 162   # P116
$ws.Cells.Item(116, 17).Value = 60   # Q116

$ws.Cells.Item(117, 4).Value = 44552    # D117
$ws.Cells.Item(117, 10).Value = 540   # J117
$ws.Cells.Item(117, 11).Value = 10500   # K117
$ws.Cells.Item(117, 12).Value = 11000   # L117
$ws.Cells.Item(117, 13).Value = 10750   # M117
$ws.Cells.Item(117, 14).Value = "`$/caja 60 unidades"   # N117
$ws.Cells.Item(117, 16).Value = 179   # P117
$ws.Cells.Item(117, 17).Value = 60   # Q117

$ws.Cells.Item(118, 4).Value = 44557    # D118
$ws.Cells.Item(118, 10).Value = 500   # J118
$ws.Cells.Item(118, 11).Value = 9500   # K118
$ws.Cells.Item(118, 12).Value = 10000   # L118
$ws.Cells.Item(118, 13).Value = 9750   # M118
$ws.Cells.Item(118, 14).Value = "`$/caja 60 unidades"   # N118
$ws.Cells.Item(118, 16).Value = 162   # P118
$ws.Cells.Item(118, 17).Value = 60   # Q118

$ws.Cells.Item(119, 4).Value = 44711    # D119
$ws.Cells.Item(119, 10).Value = 600   # J119
$ws.Cells.Item(119, 11).Value = 8000   # K119
$ws.Cells.Item(119, 12).Value = 9000   # L119
$ws.Cells.Item(119, 13).Value = 8500   # M119
$ws.Cells.Item(119, 14).Value = "`$/caja 50 unidades"   # N119
$ws.Cells.Item(119, 16).Value = 170   # P119
$ws.Cells.Item(119, 17).Value = 50   # Q119

$ws.Cells.Item(120, 4).Value = 44260    # D120
$ws.Cells.Item(120, 10).Value = 800   # J120
$ws.Cells.Item(120, 11).Value = 8000   # K120
$ws.Cells.Item(120, 12).Value = 8500   # L120
$ws.Cells.Item(120, 13).Value = 8250   # M120
$ws.Cells.Item(120, 14).Value = "`$/caja 60 unidades"   # N120
$ws.Cells.Item(120, 16).Value = 138   # P120
$ws.Cells.Item(120, 17).Value = 60   # Q120

$ws.Cells.Item(121, 4).Value = 44237    # D121
$ws.Cells.Item(121, 10).Value = 600   # J121
$ws.Cells.Item(121, 11).Value = 8000   # K121
$ws.Cells.Item(121, 12).Value = 9000   # L121
$ws.Cells.Item(121, 13).Value = 8500   # M121
$ws.Cells.Item(121, 14).Value = "`$/caja 60 unidades"   # N121
$ws.Cells.Item(121, 16).Value = 142   # P121
$ws.Cells.Item(121, 17).Value = 60   # Q121

$ws.Cells.Item(122, 4).Value = 44279    # D122
$ws.Cells.Item(122, 10).Value = 600   # J122
$ws.Cells.Item(122, 11).Value = 8000   # K122
$ws.Cells.Item(122, 12).Value = 9000   # L122
$ws.Cells.Item(122, 13).Value = 8500   # M122
$ws.Cells.Item(122, 14).Value = "`$/caja 60 unidades"   # N122
$ws.Cells.Item(122, 16).Value = 142   # P122
$ws.Cells.Item(122, 17).Value = 60   # Q122

$ws.Cells.Item(123, 4).Value = 44761    # D123
$ws.Cells.Item(123, 10).Value = 500   # J123
$ws.Cells.Item(123, 11).Value = 10000   # K123
$ws.Cells.Item(123, 12).Value = 11000   # L123
$ws.Cells.Item(123, 13).Value = 10500   # M123
$ws.Cells.Item(123, 14).Value = "`$/caja 50 unidades"   # N123
$ws.Cells.Item(123, 16).Value = 210   # P123
$ws.Cells.Item(123, 17).Value = 50   # Q123

$ws.Cells.Item(124, 4).Value = 44342    # D124
$ws.Cells.Item(124, 10).Value = 600   # J124
$ws.Cells.Item(124, 11).Value = 12000   # K124
$ws.Cells.Item(124, 12).Value = 13000   # L124
$ws.Cells.Item(124, 13).Value = 12500   # M124
$ws.Cells.Item(124, 14).Value = "`$/caja 60 unidades"   # N124
$ws.Cells.Item(124, 16).Value = 208   # P124
$ws.Cells.Item(124, 17).Value = 60   # Q124

$ws.Cells.Item(125, 4).Value = 44449    # D125
$ws.Cells.Item(125, 10).Value = 600   # J125
$ws.Cells.Item(125, 11).Value = 9000   # K125
$ws.Cells.Item(125, 12).Value = 10000   # L125
$ws.Cells.Item(125, 13).Value = 9500   # M125
$ws.Cells.Item(125, 14).Value = "`$/caja 60 unidades"   # N125
$ws.Cells.Item(125, 16).Value = 158   # P125
$ws.Cells.Item(125, 17).Value = 60   # Q125

$ws.Cells.Item(126, 4).Value = 44376    # D126
$ws.Cells.Item(126, 10).Value = 520   # J126
$ws.Cells.Item(126, 11).Value = 12000   # K126
$ws.Cells.Item(126, 12).Value = 13000   # L126
$ws.Cells.Item(126, 13).Value = 12500   # M126
$ws.Cells.Item(126, 14).Value = "`$/caja 60 unidades"   # N126
$ws.Cells.Item(126, 16).Value = 208   # P126
$ws.Cells.Item(126, 17).Value = 60   # Q126

$ws.Cells.Item(127, 4).Value = 44391    # D127
$ws.Cells.Item(127, 10).Value = 600   # J127
$ws.Cells.Item(127, 11).Value = 12000   # K127
$ws.Cells.Item(127, 12).Value = 13000   # L127
$ws.Cells.Item(127, 13).Value = 12500   # M127
$ws.Cells.Item(127, 14).Value = "`$/caja 60 unidades"   # N127
$ws.Cells.Item(127, 16).Value = 208   # P127
$ws.Cells.Item(127, 17).Value = 60   # Q127

$ws.Cells.Item(128, 4).Value = 44600    # D128
$ws.Cells.Item(128, 10).Value = 520   # J128
$ws.Cells.Item(128, 11).Value = 8500   # K128
$ws.Cells.Item(128, 12).Value = 9000   # L128
$ws.Cells.Item(128, 13).Value = 8750   # M128
$ws.Cells.Item(128, 14).Value = "`$/caja 50 unidades"   # N128
$ws.Cells.Item(128, 16).Value = 175   # P128
$ws.Cells.Item(128, 17).Value = 50   # Q128

$ws.Cells.Item(129, 4).Value = 44763    # D129
$ws.Cells.Item(129, 10).Value = 480   # J129
$ws.Cells.Item(129, 11).Value = 10000   # K129
$ws.Cells.Item(129, 12).Value = 11000   # L129
$ws.Cells.Item(129, 13).Value = 10500   # M129
$ws.Cells.Item(129, 14).Value = "`$/caja 50 unidades"   # N129
$ws.Cells.Item(129, 16).Value = 210   # P129
$ws.Cells.Item(129, 17).Value = 50   # Q129

$ws.Cells.Item(130, 4).Value = 44371    # D130
$ws.Cells.Item(130, 10).Value = 560   # J130
$ws.Cells.Item(130, 11).Value = 13000   # K130
$ws.Cells.Item(130, 12).Value = 14000   # L130
$ws.Cells.Item(130, 13).Value = 13500   # M130
$ws.Cells.Item(130, 14).Value = "`$/caja 60 unidades"   # N130
$ws.Cells.Item(130, 16).Value = 225   # P130
$ws.Cells.Item(130, 17).Value = 60   # Q130

$ws.Cells.Item(131, 4).Value = 44336    # D131
$ws.Cells.Item(131, 10).Value = 450   # J131
$ws.Cells.Item(131, 11).Value = 12000   # K131
$ws.Cells.Item(131, 12).Value = 12500   # L131
$ws.Cells.Item(131, 13).Value = 12250   # M131
$ws.Cells.Item(131, 14).Value = "`$/caja 60 unidades"   # N131
$ws.Cells.Item(131, 16).Value = 204   # P131
$ws.Cells.Item(131, 17).Value = 60   # Q131

$ws.Cells.Item(132, 4).Value = 44343    # D132
$ws.Cells.Item(132, 10).Value = 500   # J132
$ws.Cells.Item(132, 11).Value = 12000   # K132
$ws.Cells.Item(132, 12).Value = 13000   # L132
$ws.Cells.Item(132, 13).Value = 12500   # M132
$ws.Cells.Item(132, 14).Value = "`$/caja 60 unidades"   # N132
$ws.Cells.Item(132, 16).Value = 208   # P132
$ws.Cells.Item(132, 17).Value = 60   # Q132

$ws.Cells.Item(133, 4).Value = 44365    # D133
$ws.Cells.Item(133, 10).Value = 520   # J133
$ws.Cells.Item(133, 11).Value = 13000   # K133
$ws.Cells.Item(133, 12).Value = 14000   # L133
$ws.Cells.Item(133, 13).Value = 13500   # M133
$ws.Cells.Item(133, 14).Value = "`$/caja 60 unidades"   # N133
$ws.Cells.Item(133, 16).Value = 225   # P133
$ws.Cells.Item(133, 17).Value = 60   # Q133

$ws.Cells.Item(134, 4).Value = 44454    # D134
$ws.Cells.Item(134, 10).Value = 600   # J134
$ws.Cells.Item(134, 11).Value = 9000   # K134
$ws.Cells.Item(134, 12).Value = 10000   # L134
$ws.Cells.Item(134, 13).Value = 9500   # M134
$ws.Cells.Item(134, 14).Value = "`$/caja 60 unidades"   # N134
$ws.Cells.Item(134, 16).Value = 158   # P134
$ws.Cells.Item(134, 17).Value = 60   # Q134

$ws.Cells.Item(135, 4).Value = 44561    # D135
$ws.Cells.Item(135, 10).Value = 520   # J135
$ws.Cells.Item(135, 11).Value = 9000   # K135
$ws.Cells.Item(135, 12).Value = 10000   # L135
$ws.Cells.Item(135, 13).Value = 9500   # M135
$ws.Cells.Item(135, 14).Value = "`$/caja 60 unidades"   # N135
$ws.Cells.Item(135, 16).Value = 158   # P135
$ws.Cells.Item(135, 17).Value = 60   # Q135

$ws.Cells.Item(136, 4).Value = 44421    # D136
$ws.Cells.Item(136, 10).Value = 600   # J136
$ws.Cells.Item(136, 11).Value = 12000   # K136
$ws.Cells.Item(136, 12).Value = 12500   # L136
$ws.Cells.Item(136, 13).Value = 12250   # M136
$ws.Cells.Item(136, 14).Value = "`$/caja 60 unidades"   # N136
$ws.Cells.Item(136, 16).Value = 204   # P136
$ws.Cells.Item(136, 17).Value = 60   # Q136

$ws.Cells.Item(137, 4).Value = 44489    # D137
$ws.Cells.Item(137, 10).Value = 500   # J137
$ws.Cells.Item(137, 11).Value = 8000   # K137
$ws.Cells.Item(137, 12).Value = 9000   # L137
$ws.Cells.Item(137, 13).Value = 8500   # M137
$ws.Cells.Item(137, 14).Value = "`$/caja 60 unidades"   # N137
$ws.Cells.Item(137, 16).Value = 142   # P137
$ws.Cells.Item(137, 17).Value = 60   # Q137

$ws.Cells.Item(138, 4).Value = 44405    # D138
$ws.Cells.Item(138, 10).Value = 600   # J138
$ws.Cells.Item(138, 11).Value = 12000   # K138
$ws.Cells.Item(138, 12).Value = 12500   # L138
$ws.Cells.Item(138, 13).Value = 12250   # M138
$ws.Cells.Item(138, 14).Value = "`$/caja 60 unidades"   # N138
$ws.Cells.Item(138, 16).Value = 204   # P138
$ws.Cells.Item(138, 17).Value = 60   # Q138

$ws.Cells.Item(139, 4).Value = 44589    # D139
$ws.Cells.Item(139, 10).Value = 500   # J139
$ws.Cells.Item(139, 11).Value = 8000   # K139
$ws.Cells.Item(139, 12).Value = 9000   # L139
$ws.Cells.Item(139, 13).Value = 8500   # M139
$ws.Cells.Item(139, 14).Value = "`$/caja 50 unidades"   # N139
$ws.Cells.Item(139, 16).Value = 170   # P139
$ws.Cells.Item(139, 17).Value = 50   # Q139

$ws.Cells.Item(140, 4).Value = 44267    # D140
$ws.Cells.Item(140, 10).Value = 600   # J140
$ws.Cells.Item(140, 11).Value = 8000   # K140
$ws.Cells.Item(140, 12).Value = 8500   # L140
$ws.Cells.Item(140, 13).Value = 8250   # M140
$ws.Cells.Item(140, 14).Value = "`$/caja 60 unidades"   # N140
$ws.Cells.Item(140, 16).Value = 138   # P140
$ws.Cells.Item(140, 17).Value = 60   # Q140

$ws.Cells.Item(141, 4).Value = 44413    # D141
$ws.Cells.Item(141, 10).Value = 640   # J141
$ws.Cells.Item(141, 11).Value = 12000   # K141
$ws.Cells.Item(141, 12).Value = 13000   # L141
$ws.Cells.Item(141, 13).Value = 12500   # M141
$ws.Cells.Item(141, 14).Value = "`$/caja 60 unidades"   # N141
$ws.Cells.Item(141, 16).Value = 208   # P141
$ws.Cells.Item(141, 17).Value = 60   # Q141

$ws.Cells.Item(142, 4).Value = 44328    # D142
$ws.Cells.Item(142, 10).Value = 600   # J142
$ws.Cells.Item(142, 11).Value = 12000   # K142
$ws.Cells.Item(142, 12).Value = 13000   # L142
$ws.Cells.Item(142, 13).Value = 12500   # M142
$ws.Cells.Item(142, 14).Value = "`$/caja 60 unidades"   # N142
$ws.Cells.Item(142, 16).Value = 208   # P142
$ws.Cells.Item(142, 17).Value = 60   # Q142

$ws.Cells.Item(143, 4).Value = 44356    # D143
$ws.Cells.Item(143, 10).Value = 600   # J143
$ws.Cells.Item(143, 11).Value = 12000   # K143
$ws.Cells.Item(143, 12).Value = 13000   # L143
$ws.Cells.Item(143, 13).Value = 12500   # M143
$ws.Cells.Item(143, 14).Value = "`$/caja 60 unidades"   # N143
$ws.Cells.Item(143, 16).Value = 208   # P143
$ws.Cells.Item(143, 17).Value = 60   # Q143

$ws.Cells.Item(144, 4).Value = 44379    # D144
$ws.Cells.Item(144, 10).Value = 600   # J144
$ws.Cells.Item(144, 11).Value = 12000   # K144
$ws.Cells.Item(144, 12).Value = 13000   # L144
$ws.Cells.Item(144, 13).Value = 12500   # M144
$ws.Cells.Item(144, 14).Value = "`$/caja 60 unidades"   # N144
$ws.Cells.Item(144, 16).Value = 208   # P144
$ws.Cells.Item(144, 17).Value = 60   # Q144

$ws.Cells.Item(145, 4).Value = 44322    # D145
$ws.Cells.Item(145, 10).Value = 440   # J145
$ws.Cells.Item(145, 11).Value = 8000   # K145
$ws.Cells.Item(145, 12).Value = 9000   # L145
$ws.Cells.Item(145, 13).Value = 8500   # M145
$ws.Cells.Item(145, 14).Value = "`$/caja 60 unidades"   # N145
$ws.Cells.Item(145, 16).Value = 142   # P145
$ws.Cells.Item(145, 17).Value = 60   # Q145

$ws.Cells.Item(146, 4).Value = 44782    # D146
$ws.Cells.Item(146, 10).Value = 540   # J146
$ws.Cells.Item(146, 11).Value = 10000   # K146
$ws.Cells.Item(146, 12).Value = 11000   # L146
$ws.Cells.Item(146, 13).Value = 10500   # M146
$ws.Cells.Item(146, 14).Value = "`$/caja 50 unidades"   # N146
$ws.Cells.Item(146, 16).Value = 210   # P146
$ws.Cells.Item(146, 17).Value = 50   # Q146

$ws.Cells.Item(147, 4).Value = 44251    # D147
$ws.Cells.Item(147, 10).Value = 600   # J147
$ws.Cells.Item(147, 11).Value = 8000   # K147
$ws.Cells.Item(147, 12).Value = 9000   # L147
$ws.Cells.Item(147, 13).Value = 8500   # M147
$ws.Cells.Item(147, 14).Value = "`$/caja 60 unidades"   # N147
$ws.Cells.Item(147, 16).Value = 142   # P147
$ws.Cells.Item(147, 17).Value = 60   # Q147

$ws.Cells.Item(148, 4).Value = 44344    # D148
$ws.Cells.Item(148, 10).Value = 520   # J148
$ws.Cells.Item(148, 11).Value = 12000   # K148
$ws.Cells.Item(148, 12).Value = 13000   # L148
$ws.Cells.Item(148, 13).Value = 12500   # M148
$ws.Cells.Item(148, 14).Value = "`$/caja 60 unidades"   # N148
$ws.Cells.Item(148, 16).Value = 208   # P148
$ws.Cells.Item(148, 17).Value = 60   # Q148

$ws.Cells.Item(149, 4).Value = 44455    # D149
$ws.Cells.Item(149, 10).Value = 580   # J149
$ws.Cells.Item(149, 11).Value = 9000   # K149
$ws.Cells.Item(149, 12).Value = 10000   # L149
$ws.Cells.Item(149, 13).Value = 9500   # M149
$ws.Cells.Item(149, 14).Value = "`$/caja 60 unidades"   # N149
$ws.Cells.Item(149, 16).Value = 158   # P149
$ws.Cells.Item(149, 17).Value = 60   # Q149

$ws.Cells.Item(150, 4).Value = 44484    # D150
$ws.Cells.Item(150, 10).Value = 600   # J150
$ws.Cells.Item(150, 11).Value = 9000   # K150
$ws.Cells.Item(150, 12).Value = 10000   # L150
$ws.Cells.Item(150, 13).Value = 9500   # M150
$ws.Cells.Item(150, 14).Value = "`$/caja 60 unidades"   # N150
$ws.Cells.Item(150, 16).Value = 158   # P150
$ws.Cells.Item(150, 17).Value = 60   # Q150

$ws.Cells.Item(151, 4).Value = 44665    # D151
$ws.Cells.Item(151, 10).Value = 400   # J151
$ws.Cells.Item(151, 11).Value = 8000   # K151
$ws.Cells.Item(151, 12).Value = 9000   # L151
$ws.Cells.Item(151, 13).Value = 8500   # M151
$ws.Cells.Item(151, 14).Value = "`$/caja 50 unidades"   # N151
$ws.Cells.Item(151, 16).Value = 170   # P151
$ws.Cells.Item(151, 17).Value = 50   # Q151

$ws.Cells.Item(152, 4).Value = 44452    # D152
$ws.Cells.Item(152, 10).Value = 560   # J152
$ws.Cells.Item(152, 11).Value = 9000   # K152
$ws.Cells.Item(152, 12).Value = 10000   # L152
$ws.Cells.Item(152, 13).Value = 9500   # M152
$ws.Cells.Item(152, 14).Value = "`$/caja 60 unidades"   # N152
$ws.Cells.Item(152, 16).Value = 158   # P152
$ws.Cells.Item(152, 17).Value = 60   # Q152

$ws.Cells.Item(153, 4).Value = 44510    # D153
$ws.Cells.Item(153, 10).Value = 520   # J153
$ws.Cells.Item(153, 11).Value = 8000   # K153
$ws.Cells.Item(153, 12).Value = 8500   # L153
$ws.Cells.Item(153, 13).Value = 8250   # M153
$ws.Cells.Item(153, 14).Value = "`$/caja 60 unidades"   # N153
$ws.Cells.Item(153, 16).Value = 138   # P153
$ws.Cells.Item(153, 17).Value = 60   # Q153

$ws.Cells.Item(154, 4).Value = 44516    # D154
$ws.Cells.Item(154, 10).Value = 400   # J154
$ws.Cells.Item(154, 11).Value = 8000   # K154
$ws.Cells.Item(154, 12).Value = 9000   # L154
$ws.Cells.Item(154, 13).Value = 8500   # M154
$ws.Cells.Item(154, 14).Value = "`$/caja 60 unidades"   # N154
$ws.Cells.Item(154, 16).Value = 142   # P154
$ws.Cells.Item(154, 17).Value = 60   # Q154

$ws.Cells.Item(155, 4).Value = 44463    # D155
$ws.Cells.Item(155, 10).Value = 600   # J155
$ws.Cells.Item(155, 11).Value = 9000   # K155
$ws.Cells.Item(155, 12).Value = 10000   # L155
$ws.Cells.Item(155, 13).Value = 9500   # M155
$ws.Cells.Item(155, 14).Value = "`$/caja 60 unidades"   # N155
$ws.Cells.Item(155, 16).Value = 158   # P155
$ws.Cells.Item(155, 17).Value = 60   # Q155

$ws.Cells.Item(156, 4).Value = 44382    # D156
$ws.Cells.Item(156, 10).Value = 560   # J156
$ws.Cells.Item(156, 11).Value = 12000   # K156
$ws.Cells.Item(156, 12).Value = 13000   # L156
$ws.Cells.Item(156, 13).Value = 12500   # M156
$ws.Cells.Item(156, 14).Value = "`$/caja 60 unidades"   # N156
$ws.Cells.Item(156, 16).Value = 208   # P156
$ws.Cells.Item(156, 17).Value = 60   # Q156
